$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31 (shifts existing rows 31-105 down to 32-106)
$ws.Rows(31).Insert()

# Populate the new row 31 with the new record
$ws.Cells.Item(31, 1).Value = 11
$ws.Cells.Item(31, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(31, 3).Value = 'Bíobío'
$ws.Cells.Item(31, 4).Value = 44622
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(31, 6).Value = 100112043
$ws.Cells.Item(31, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(31, 8).Value = 'Sin especificar'
$ws.Cells.Item(31, 9).Value = 'Primera'
$ws.Cells.Item(31, 10).Value = 250
$ws.Cells.Item(31, 11).Value = 8500
$ws.Cells.Item(31, 12).Value = 9000
$ws.Cells.Item(31, 13).Value = 8800
$ws.Cells.Item(31, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(31, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(31, 16).Value = 147
$ws.Cells.Item(31, 17).Value = 60
$ws.Cells.Item(31, 18).Value = 'Hortaliza'
